$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row labels:
#    "<name>_old" -> "<name>_FV2310" (columns A-J)
#    "<name>_new" -> "<name>_FV2404" (columns L-U)
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2310')
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2404')
}

# 2. Freeze the header row (split above row 2 / A2) in the sheet view.
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the data range into an Excel Table ("Table1") spanning A1:U72.
$tableRange = $ws.Range("A1:U72")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
